$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 106, shifting existing rows 106-217 down to 107-218.
$ws.Rows.Item(106).Insert()

# Populate the newly inserted row 106 with its data.
$ws.Cells.Item(106, 1).Value = 10
$ws.Cells.Item(106, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(106, 3).Value = "La Araucanía"
$ws.Cells.Item(106, 4).Value = 44904
$ws.Cells.Item(106, 5).Value = 9
$ws.Cells.Item(106, 6).Value = 100112012
$ws.Cells.Item(106, 7).Value = "Espinaca"
$ws.Cells.Item(106, 8).Value = "Sin especificar"
$ws.Cells.Item(106, 9).Value = "Primera"
$ws.Cells.Item(106, 10).Value = 60
$ws.Cells.Item(106, 11).Value = 9000
$ws.Cells.Item(106, 12).Value = 10000
$ws.Cells.Item(106, 13).Value = 9417
$ws.Cells.Item(106, 14).Value = "$/docena de atados"
$ws.Cells.Item(106, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(106, 16).Value = 3139
$ws.Cells.Item(106, 17).Value = 3
$ws.Cells.Item(106, 18).Value = "Hortaliza"
